# Apply MM_data.xlsx update: renumber per-city SortOrder, fix row 101
# (Awaali Gardens -> Etrah Garden), add 3 new MADINAH/food rows, and
# refresh UI/view metadata (selection, autofilter range, filter defined name).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row 101 content change: "Awaali Gardens" (+ details) is replaced by
#    a brand-new place "Etrah Garden" with no Details text.
# ---------------------------------------------------------------------
$ws.Range("D101").Value = "Etrah Garden"
$ws.Range("E101").ClearContents()

# ---------------------------------------------------------------------
# 2) Re-number the SortOrder (column K) values for every existing data
#    row (2-105) per the new per-city ordering.
# ---------------------------------------------------------------------
$ws.Range("K2").Value = 18
$ws.Range("K3").Value = 16
$ws.Range("K4").Value = 19
$ws.Range("K5").Value = 17
$ws.Range("K6").Value = 26
$ws.Range("K7").Value = 5
$ws.Range("K8").Value = 31
$ws.Range("K9").Value = 32
$ws.Range("K10").Value = 22
$ws.Range("K11").Value = 46
$ws.Range("K12").Value = 20
$ws.Range("K13").Value = 27
$ws.Range("K14").Value = 6
$ws.Range("K15").Value = 23
$ws.Range("K16").Value = 33
$ws.Range("K17").Value = 2
$ws.Range("K18").Value = 18
$ws.Range("K19").Value = 11
$ws.Range("K20").Value = 39
$ws.Range("K21").Value = 12
$ws.Range("K22").Value = 34
$ws.Range("K23").Value = 1
$ws.Range("K24").Value = 28
$ws.Range("K25").Value = 7
$ws.Range("K26").Value = 40
$ws.Range("K28").Value = 35
$ws.Range("K29").Value = 13
$ws.Range("K30").Value = 36
$ws.Range("K31").Value = 7
$ws.Range("K32").Value = 12
$ws.Range("K33").Value = 29
$ws.Range("K34").Value = 21
$ws.Range("K35").Value = 22
$ws.Range("K36").Value = 30
$ws.Range("K37").Value = 19
$ws.Range("K38").Value = 1
$ws.Range("K39").Value = 37
$ws.Range("K40").Value = 38
$ws.Range("K41").Value = 39
$ws.Range("K42").Value = 40
$ws.Range("K43").Value = 41
$ws.Range("K44").Value = 42
$ws.Range("K45").Value = 31
$ws.Range("K46").Value = 32
$ws.Range("K47").Value = 3
$ws.Range("K48").Value = 44
$ws.Range("K49").Value = 33
$ws.Range("K50").Value = 34
$ws.Range("K51").Value = 23
$ws.Range("K52").Value = 24
$ws.Range("K53").Value = 1
$ws.Range("K54").Value = 13
$ws.Range("K55").Value = 3
$ws.Range("K56").Value = 43
$ws.Range("K57").Value = 44
$ws.Range("K58").Value = 8
$ws.Range("K59").Value = 24
$ws.Range("K60").Value = 45
$ws.Range("K61").Value = 14
$ws.Range("K62").Value = 9
$ws.Range("K63").Value = 45
$ws.Range("K64").Value = 14
$ws.Range("K65").Value = 15
$ws.Range("K66").Value = 25
$ws.Range("K67").Value = 46
$ws.Range("K68").Value = 16
$ws.Range("K69").Value = 47
$ws.Range("K70").Value = 9
$ws.Range("K71").Value = 4
$ws.Range("K72").Value = 10
$ws.Range("K73").Value = 11
$ws.Range("K74").Value = 15
$ws.Range("K75").Value = 35
$ws.Range("K77").Value = 26
$ws.Range("K78").Value = 27
$ws.Range("K79").Value = 28
$ws.Range("K80").Value = 42
$ws.Range("K81").Value = 20
$ws.Range("K82").Value = 48
$ws.Range("K83").Value = 49
$ws.Range("K84").Value = 50
$ws.Range("K85").Value = 10
$ws.Range("K86").Value = 51
$ws.Range("K87").Value = 5
$ws.Range("K88").Value = 25
$ws.Range("K89").Value = 21
$ws.Range("K90").Value = 8
$ws.Range("K91").Value = 17
$ws.Range("K92").Value = 52
$ws.Range("K93").Value = 53
$ws.Range("K94").Value = 29
$ws.Range("K95").Value = 54
$ws.Range("K96").Value = 37
$ws.Range("K97").Value = 2
$ws.Range("K98").Value = 43
$ws.Range("K99").Value = 30
$ws.Range("K100").Value = 55
$ws.Range("K101").Value = 48
$ws.Range("K102").Value = 6
$ws.Range("K103").Value = 38
$ws.Range("K104").Value = 56
$ws.Range("K105").Value = 47

# ---------------------------------------------------------------------
# 3) Append three new MADINAH / food rows (107-109). The Name values are
#    written bottom-row-first (109, 108, 107) to reproduce the original
#    shared-string insertion order (Ajwatech, Saraya Ward Restaurant,
#    Reef al-Maknan).
# ---------------------------------------------------------------------
$ws.Range("D109").Value = "Ajwatech"
$ws.Range("D108").Value = "Saraya Ward Restaurant"
$ws.Range("D107").Value = "Reef al-Maknan"

$ws.Range("A107").Value = "MADINAH"
$ws.Range("B107").Value = "food"
$ws.Range("G107").Value = "pending"
$ws.Range("K107").Value = 50
$ws.Range("L107").Value = "No"
$ws.Range("M107").Value = "No"

$ws.Range("A108").Value = "MADINAH"
$ws.Range("B108").Value = "food"
$ws.Range("G108").Value = "pending"
$ws.Range("K108").Value = 49
$ws.Range("L108").Value = "No"
$ws.Range("M108").Value = "No"

$ws.Range("A109").Value = "MADINAH"
$ws.Range("B109").Value = "food"
$ws.Range("G109").Value = "pending"
$ws.Range("K109").Value = 51
$ws.Range("L109").Value = "No"
$ws.Range("M109").Value = "No"

# ---------------------------------------------------------------------
# 4) Refresh the AutoFilter range to cover the newly added rows (toggle
#    off first so re-applying doesn't just strip filtering instead).
# ---------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:M109").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name in sync with the
# refreshed AutoFilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$M`$109"
    }
}

# ---------------------------------------------------------------------
# 5) Update the active selection to match the saved view state.
# ---------------------------------------------------------------------
$ws.Range("S114").Select()
